$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projections")

$ws.Range("A4").Value = 4831
$ws.Range("B4").Value = 733
$ws.Range("C4").Value = 4268
$ws.Range("D4").Value = 2.195954545454544
$ws.Range("E4").Value = 0.3333636363636364
$ws.Range("F4").Value = 1.940409090909091
$ws.Range("G4").Value = 0.3838636363636364
$ws.Range("H4").Value = 0.08254545454545451
$ws.Range("I4").Value = 0.291681818181818
$ws.Range("J4").Value = 0.2013125
$ws.Range("K4").Value = 0.1556416666666667
$ws.Range("L4").Value = 0.1935491666666667
$ws.Range("M4").Value = 3021
$ws.Range("N4").Value = 136
$ws.Range("O4").Value = 7456
$ws.Range("P4").Value = 1.756744186046511
$ws.Range("Q4").Value = 0.07965116279069764
$ws.Range("R4").Value = 4.335174418604652
$ws.Range("S4").Value = 0.3188372093023256
$ws.Range("T4").Value = 0.02017441860465117
$ws.Range("U4").Value = 0.6362790697674418
$ws.Range("V4").Value = 0.2375275
$ws.Range("W4").Value = 0.1572375
$ws.Range("X4").Value = 0.2486216666666666
